$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 63, shifting existing rows 63-104 down to 64-105.
$ws.Rows("63:63").Insert()

# Populate the newly inserted row 63 with the new data point.
$ws.Range("A63").Value = 8
$ws.Range("B63").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C63").Value = 'Coquimbo'
$ws.Range("D63").Value = 44767
$ws.Range("E63").Value = 4
$ws.Range("F63").Value = 100112052
$ws.Range("G63").Value = 'Albahaca'
$ws.Range("H63").Value = 'Sin especificar'
$ws.Range("I63").Value = 'Primera'
$ws.Range("J63").Value = 1200
$ws.Range("K63").Value = 3500
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = 3750
$ws.Range("N63").Value = '$/paquete'
$ws.Range("O63").Value = 'Región de Arica y Parinacota'
$ws.Range("P63").Value = 3750
$ws.Range("Q63").Value = 1
$ws.Range("R63").Value = 'Hortaliza'
